# Test Cases updated for Header Scenarios
# Populate Sheet1 with the TestCase header row + first sample row, size the
# two used columns to fit their content, and leave the selection where the
# author left it (B5) like the saved workbook shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "TestCase ID"
$ws.Range("B1").Value = "TestCaseName"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Header validation"

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

$ws.Range("B5").Select()
